$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '43.665.86'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +4.17%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.261.38'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +1.24%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '230.69'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.41%  '
$ws.Range('E6').Value = '  -0.07%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '61.27'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -0.98%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').Value = '  +4.47%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '58.01'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -2.39%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0934'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +4.95%  '
$ws.Range('E12').Value = '  +0.51%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '2.600.09'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +1.29%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '15.60'
$ws.Range('D14').Style = "Normal"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '23.60'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +7.24%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '5.80'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +3.59%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.809'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +0.94%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.252.73'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.34%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '43.283.80'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +3.58%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.0₃0935'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +4.22%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '72.90'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.92%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.21'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +2.55%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '253.04'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +1.25%  '
$ws.Range('E24').Value = '  -0.09%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.54'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +6.15%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.37'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.83'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +1.48%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '170.69'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +2.40%  '
$ws.Range('E29').Value = '  -1.41%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '20.47'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +2.52%  '
$ws.Range('E31').Value = '  +1.70%  '
$ws.Range('E33').Value = '  -0.12%  '
$ws.Range('E34').Value = '  +1.40%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.79'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +2.16%  '
$ws.Range('E36').Value = '  +3.83%  '
$ws.Range('E37').Value = '  -2.91%  '
$ws.Range('E38').Value = '  +1.21%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.60'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -1.51%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0250'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +4.10%  '
$ws.Range('E41').Value = '  -0.08%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.000229'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -10.04%  '
$ws.Range('E43').Value = '  +1.46%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0994'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +1.46%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '4.51'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -7.06%  '
$ws.Range('E46').Value = '  -0.51%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '98.18'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.75%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.471.82'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.56%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '16.65'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +0.90%  '
$ws.Range('E50').Value = '  +0.49%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.26'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +7.61%  '
